# Update the "Antibodies" series data (first bar-chart series on the
# single chart on slide 1) with the revised percentages from the
# re-analysed raw data.
#
#   0.967441860465116 -> 0.964556962025316
#   0.646511627906977 -> 0.615189873417721
#   0.613953488372093 -> 0.579746835443038

$p    = $ppt.ActivePresentation
$s    = $p.Slides.Item(1)
$shp  = $s.Shapes.Item(1)
$chart = $shp.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Values = @(0.964556962025316, 0.615189873417721, 0.579746835443038)
